# Update the NATMI LR-pairs sheet (Fgf22-Fgfrl1) with freshly-recomputed TPM
# values. The sending/target cluster assignments are refreshed too: the
# first block of rows (previously "FAPs") now sends from "ECs", and the
# second block (previously "MuSCs") now sends from "FAPs" -- ligand/receptor
# symbols (Fgf22 / Fgfrl1) and target clusters (ECs / FAPs / MuSCs) are
# unchanged in meaning.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One row per data record (row 2 .. row 7 on the sheet). Column letters map
# 1:1 onto the sheet's existing header row.
$rows = @(
    @{ Row = 2;  A = "ECs";  B = "Fgf22"; C = "Fgfrl1"; D = "ECs";
       E = 1; F = 0.3333333333333333; G = 0.028883; H = 0.086649;
       I = 0.2501761214025038; J = 0.2501761214025038; K = 3; L = 1;
       M = 0.259369; N = 0.778107; O = 0.0514155333512404; P = 0.0514155333512404;
       Q = 0.007491354827000001; R = 0.067422193443;
       S = 0.0128629387136544; T = 0.0128629387136544 },

    @{ Row = 3;  A = "ECs";  B = "Fgf22"; C = "Fgfrl1"; D = "FAPs";
       E = 1; F = 0.3333333333333333; G = 0.028883; H = 0.086649;
       I = 0.2501761214025038; J = 0.2501761214025038; K = 3; L = 1;
       M = 3.018163; N = 9.054489; O = 0.5982999525231611; P = 0.5982999525231611;
       Q = 0.08717360192900001; R = 0.7845624173610001;
       S = 0.1496803615575467; T = 0.1496803615575467 },

    @{ Row = 4;  A = "ECs";  B = "Fgf22"; C = "Fgfrl1"; D = "MuSCs";
       E = 1; F = 0.3333333333333333; G = 0.028883; H = 0.086649;
       I = 0.2501761214025038; J = 0.2501761214025038; K = 3; L = 1;
       M = 1.767033; N = 5.301099; O = 0.3502845141255985; P = 0.3502845141255985;
       Q = 0.051037214139; R = 0.459334927251;
       S = 0.0876328211313028; T = 0.0876328211313028 },

    @{ Row = 5;  A = "FAPs"; B = "Fgf22"; C = "Fgfrl1"; D = "ECs";
       E = 1; F = 0.3333333333333333; G = 0.08656766666666667; H = 0.259703;
       I = 0.7498238785974961; J = 0.7498238785974962; K = 3; L = 1;
       M = 0.259369; N = 0.778107; O = 0.0514155333512404; P = 0.0514155333512404;
       Q = 0.02245296913566667; R = 0.202076722221;
       S = 0.038552594637586; T = 0.038552594637586 },

    @{ Row = 6;  A = "FAPs"; B = "Fgf22"; C = "Fgfrl1"; D = "FAPs";
       E = 1; F = 0.3333333333333333; G = 0.08656766666666667; H = 0.259703;
       I = 0.7498238785974961; J = 0.7498238785974962; K = 3; L = 1;
       M = 3.018163; N = 9.054489; O = 0.5982999525231611; P = 0.5982999525231611;
       Q = 0.2612753285296667; R = 2.351477956767;
       S = 0.4486195909656145; T = 0.4486195909656145 },

    @{ Row = 7;  A = "FAPs"; B = "Fgf22"; C = "Fgfrl1"; D = "MuSCs";
       E = 1; F = 0.3333333333333333; G = 0.08656766666666667; H = 0.259703;
       I = 0.7498238785974961; J = 0.7498238785974962; K = 3; L = 1;
       M = 1.767033; N = 5.301099; O = 0.3502845141255985; P = 0.3502845141255985;
       Q = 0.152967923733; R = 1.376711313597;
       S = 0.2626516929942957; T = 0.2626516929942957 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($rec in $rows) {
    $r = $rec.Row
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $rec[$col]
    }
}
